$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "any.admin@alyaconsulting.ch" notes for Windows 365 Administrator / Edge Administrator rows ---
$ws.Range("D17").Value = "any.admin@alyaconsulting.ch"
$ws.Range("D18").Value = "any.admin@alyaconsulting.ch"

# --- New roles appended at the bottom of the role list ---
$ws.Range("A93").Value = "Yammer Administrator"
$ws.Range("A94").Value = "Lifecycle Workflows Administrator"
$ws.Range("A95").Value = "Permissions Management Administrator"

# Rows 93-96 in column A move from the "indented/blank" style to the regular
# left-aligned text style used throughout the rest of column A.
$ws.Range("A93:A96").IndentLevel = 0

# --- Keep the sheet's remembered sort range in sync with the data (A8:A96) ---
# A helper column is used as the sort key so the existing row order (which is
# not a plain alphabetical sort) is left untouched; only the bookkeeping
# "sortState" range is refreshed.
for ($i = 8; $i -le 96; $i++) {
    $ws.Cells.Item($i, 26).Value = $i
}

$sort = $ws.Sort
Write-Output "Preparing sort: $sort"
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("Z8:Z96"))
$sort.SetRange($ws.Range("A8:A96"))
$sort.Header = 2
$sort.Apply()

$ws.Range("Z8:Z96").ClearContents()

# --- Update the active selection shown when the sheet is reopened ---
$ws.Range("C91").Select()
